# feat: add 2022-Q3 data
#
# The workbook ships one sheet per quarter ("2022-Q2" so far) plus a "总计"
# (totals) roll-up sheet. This adds a new "2022-Q3" sheet with that
# quarter's fund data, inserts it right after "总计" (so tab order becomes
# 总计, 2022-Q3, 2022-Q2), and records the new quarter in the totals sheet
# while keeping the old 2022-Q2 totals row intact.

# Writes $text into $range as a literal string, even when it looks like a
# number (e.g. "015205", "0.24") so a leading zero / exact decimal text
# survives instead of being coerced to a numeric cell - exactly like Excel
# keeps a typed value literal when it was produced by a text formula.
function Set-TextValue($range, $text) {
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy($range)
}

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(1)
$wsQ2 = $wb.Worksheets.Item(2)

# --- 1. Duplicate the existing "2022-Q2" sheet right after itself, so one
#        copy can stay untouched as "2022-Q2" while the original (position
#        2) gets repurposed into the new "2022-Q3" sheet. This keeps the
#        original sheet's sheetId (2) on the new "2022-Q3" tab and hands
#        the fresh sheetId (3) to the sheet that keeps the "2022-Q2" name -
#        matching how the workbook is ordered/identified after the edit.
$wsQ2.Copy($null, $wsQ2)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ2Copy = $wb.Worksheets.Item(3)

$wsQ3.Name = "2022-Q3"
$wsQ2Copy.Name = "2022-Q2"

# --- 2. Re-style the "2022-Q3" sheet's header row and first column to match
#        the "总计" sheet's header formatting (bold + thin border), then
#        overwrite every cell with the new quarter's fund data.
$headerStyleSrc = $wsTotal.Range("B1")
$colAStyleSrc = $wsTotal.Range("A2")

foreach ($col in @("B", "C", "D", "E", "F", "G", "H")) {
    $headerStyleSrc.Copy($wsQ3.Range($col + "1"))
}
$colAStyleSrc.Copy($wsQ3.Range("A2"))
$colAStyleSrc.Copy($wsQ3.Range("A3"))

Set-TextValue $wsQ3.Range("B1") "基金代码"
Set-TextValue $wsQ3.Range("C1") "基金名称"
Set-TextValue $wsQ3.Range("D1") "基金规模"
Set-TextValue $wsQ3.Range("E1") "股票总仓位"
Set-TextValue $wsQ3.Range("F1") "仓位占比"
Set-TextValue $wsQ3.Range("G1") "持有市值(亿元)"
Set-TextValue $wsQ3.Range("H1") "仓位排名"

$wsQ3.Range("A2").Value = 0
Set-TextValue $wsQ3.Range("B2") "015205"
Set-TextValue $wsQ3.Range("C2") "银华全球新能源车量化优选股票（QDII）C"
Set-TextValue $wsQ3.Range("D2") "0.24"
Set-TextValue $wsQ3.Range("E2") "90.37"
Set-TextValue $wsQ3.Range("F2") "7.03"
Set-TextValue $wsQ3.Range("G2") "0.0169"
$wsQ3.Range("H2").Value = 2

$wsQ3.Range("A3").Value = 1
Set-TextValue $wsQ3.Range("B3") "015204"
Set-TextValue $wsQ3.Range("C3") "银华全球新能源车量化优选股票（QDII）A"
Set-TextValue $wsQ3.Range("D3") "0.19"
Set-TextValue $wsQ3.Range("E3") "90.37"
Set-TextValue $wsQ3.Range("F3") "7.03"
Set-TextValue $wsQ3.Range("G3") "0.0134"
$wsQ3.Range("H3").Value = 2

# --- 3. "总计": row 2 now reports the new 2022-Q3 totals, and a new row 3
#        re-adds the 2022-Q2 totals that used to sit in row 2.
$wsTotal.Range("A2").Copy($wsTotal.Range("A3"))
$wsTotal.Range("A3").Value = 1
Set-TextValue $wsTotal.Range("B3") "2022-Q2"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.02

Set-TextValue $wsTotal.Range("B2") "2022-Q3"
$wsTotal.Range("D2").Value = 0.03
